$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "_old" -> "_FV2210", "_new" -> "_FV2304" ---
$oldToNew = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $oldToNew.Keys) {
    $ws.Range($addr).Value = $oldToNew[$addr]
}

# --- 2. Turn the header + data range into an Excel Table, preserving the ---
#        existing header-row formatting (bold/fill/border/center/wrap)
#        without the engine auto-capturing it into a new dxf/cellXf.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A1000:U1000")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# --- 3. Freeze the header row (split after row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
